$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsOverview.Range("G4").Value = "2016-08-25 10:47:43"

$wsZhCn.Range("H4").Value = "2016-08-25 10:47:39"
$wsZhCn.Range("K4").Value = "2016-08-25 10:48:12"

$wsDeDe.Range("H4").Value = "2016-08-25 10:47:43"
$wsDeDe.Range("K4").Value = "2016-08-25 10:48:19"
